# Auto-generated edit script: apply market-price refresh values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7250.5
$ws.Range("I40").Value = 2001
$ws.Range("K40").Value = 2001
$ws.Range("M40").Value = -1826

$ws.Range("H41").Value = 1401
$ws.Range("I41").Value = 2145.8333
$ws.Range("J41").Value = 842.375
$ws.Range("K41").Value = 2145.8333
$ws.Range("L41").Value = 842.375
$ws.Range("M41").Value = -1705.8333
$ws.Range("N41").Value = -1722.375

$ws.Range("H53").Value = 85.5
$ws.Range("I53").Value = 82.8
$ws.Range("J53").Value = 99
$ws.Range("K53").Value = 82.8
$ws.Range("L53").Value = 99
$ws.Range("M53").Value = 554.2
$ws.Range("N53").Value = -1373

$ws.Range("H76").Value = 12488.091
$ws.Range("I76").Value = 12466.667
$ws.Range("J76").Value = 12496.125
$ws.Range("K76").Value = 12466.667
$ws.Range("L76").Value = 12496.125
$ws.Range("M76").Value = -12151.667
$ws.Range("N76").Value = -13126.125

$ws.Range("H79").Value = 12488.091
$ws.Range("I79").Value = 12466.667
$ws.Range("J79").Value = 12496.125
$ws.Range("K79").Value = 12466.667
$ws.Range("L79").Value = 12496.125
$ws.Range("M79").Value = -11374.667
$ws.Range("N79").Value = -14680.125

$ws.Range("H131").Value = 1925
$ws.Range("I131").Value = 1925
$ws.Range("K131").Value = 5775
$ws.Range("M131").Value = -735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3523.6667
$ws.Range("J45").Value = 4066.6445
$ws.Range("L45").Value = 4066.6445
$ws.Range("N45").Value = -4820.6445

$ws.Range("H61").Value = 5155.033
$ws.Range("I61").Value = 3507.8333
$ws.Range("K61").Value = 3507.8333
$ws.Range("M61").Value = -3295.8333

$ws.Range("H110").Value = 1964.4546
$ws.Range("I110").Value = 1734.3334
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1734.3334
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 310.6666
$ws.Range("N110").Value = -7090

$ws.Range("H132").Value = 2123.16
$ws.Range("I132").Value = 1868.409
$ws.Range("K132").Value = 5605.227000000001
$ws.Range("M132").Value = -3075.227000000001

$ws.Range("H136").Value = 5155.033
$ws.Range("I136").Value = 3507.8333
$ws.Range("K136").Value = 10523.4999
$ws.Range("M136").Value = -7973.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 703.3125
$ws.Range("J80").Value = 835.53845
$ws.Range("L80").Value = 835.53845
$ws.Range("N80").Value = -2831.53845

$ws.Range("H83").Value = 703.3125
$ws.Range("J83").Value = 835.53845
$ws.Range("L83").Value = 4177.69225
$ws.Range("N83").Value = -14161.69225

$ws.Range("H99").Value = 4571.273
$ws.Range("I99").Value = 4571.273
$ws.Range("K99").Value = 4571.273
$ws.Range("M99").Value = -3073.273

$ws.Range("H134").Value = 3017.7273
$ws.Range("I134").Value = 3017.7273
$ws.Range("K134").Value = 9053.1819
$ws.Range("M134").Value = -6518.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 27799886
$ws.Range("I86").Value = 27799886
$ws.Range("K86").Value = 27799886
$ws.Range("M86").Value = -27798763

$ws.Range("H89").Value = 27799886
$ws.Range("I89").Value = 27799886
$ws.Range("K89").Value = 138999430
$ws.Range("M89").Value = -138993814

$ws.Range("H141").Value = 294237.38
$ws.Range("J141").Value = 294237.38
$ws.Range("L141").Value = 294237.38
$ws.Range("N141").Value = -304597.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4623.8237
$ws.Range("I56").Value = 4623.8237
$ws.Range("K56").Value = 4623.8237
$ws.Range("M56").Value = -4093.8237

$ws.Range("H114").Value = 3862.6667
$ws.Range("I114").Value = 794.5
$ws.Range("K114").Value = 2383.5
$ws.Range("M114").Value = 870.5

$ws.Range("H137").Value = 4526.7896
$ws.Range("J137").Value = 5354.1113
$ws.Range("L137").Value = 16062.3339
$ws.Range("N137").Value = -26262.3339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15366120
$ws.Range("I11").Value = 17223714
$ws.Range("K11").Value = 17223714
$ws.Range("M11").Value = -17223575

$ws.Range("H18").Value = 3343334.8
$ws.Range("I18").Value = 20005
$ws.Range("J18").Value = 5004999.5
$ws.Range("K18").Value = 20005
$ws.Range("L18").Value = 5004999.5
$ws.Range("M18").Value = -19712
$ws.Range("N18").Value = -5005585.5

$ws.Range("H23").Value = 497.5
$ws.Range("J23").Value = 495
$ws.Range("L23").Value = 495
$ws.Range("N23").Value = -941

$ws.Range("H24").Value = 1016400
$ws.Range("I24").Value = 2504000
$ws.Range("J24").Value = 24666.666
$ws.Range("K24").Value = 2504000
$ws.Range("L24").Value = 24666.666
$ws.Range("M24").Value = -2503827
$ws.Range("N24").Value = -25012.666

$ws.Range("H25").Value = 2000
$ws.Range("J25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = -3058

$ws.Range("H70").Value = 7981.909
$ws.Range("I70").Value = 7937.625
$ws.Range("J70").Value = 8100
$ws.Range("K70").Value = 7937.625
$ws.Range("L70").Value = 8100
$ws.Range("M70").Value = -7667.625
$ws.Range("N70").Value = -8640

$ws.Range("H73").Value = 7981.909
$ws.Range("I73").Value = 7937.625
$ws.Range("J73").Value = 8100
$ws.Range("K73").Value = 7937.625
$ws.Range("L73").Value = 8100
$ws.Range("M73").Value = -7001.625
$ws.Range("N73").Value = -9972

$ws.Range("H80").Value = 5835
$ws.Range("I80").Value = 4515.9165
$ws.Range("J80").Value = 8473.166999999999
$ws.Range("K80").Value = 4515.9165
$ws.Range("L80").Value = 8473.166999999999
$ws.Range("M80").Value = -3517.9165
$ws.Range("N80").Value = -10469.167

$ws.Range("H83").Value = 5835
$ws.Range("I83").Value = 4515.9165
$ws.Range("J83").Value = 8473.166999999999
$ws.Range("K83").Value = 22579.5825
$ws.Range("L83").Value = 42365.835
$ws.Range("M83").Value = -17587.5825
$ws.Range("N83").Value = -52349.835

$ws.Range("H97").Value = 839.75
$ws.Range("I97").Value = 839.75
$ws.Range("K97").Value = 839.75
$ws.Range("M97").Value = -343.75

$ws.Range("H99").Value = 23239.834
$ws.Range("I99").Value = 13888
$ws.Range("K99").Value = 13888
$ws.Range("M99").Value = -11642

$ws.Range("H111").Value = 44195
$ws.Range("J111").Value = 44195
$ws.Range("L111").Value = 44195
$ws.Range("N111").Value = -50329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 9499.799999999999
$ws.Range("J10").Value = 4500
$ws.Range("L10").Value = 4500
$ws.Range("N10").Value = -4780

$ws.Range("H16").Value = 9932.333000000001
$ws.Range("I16").Value = 9899
$ws.Range("K16").Value = 9899
$ws.Range("M16").Value = -9729

$ws.Range("H20").Value = 2516000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20452

$ws.Range("H46").Value = 3712.25
$ws.Range("I46").Value = 4999
$ws.Range("J46").Value = 3283.3333
$ws.Range("K46").Value = 4999
$ws.Range("L46").Value = 3283.3333
$ws.Range("M46").Value = -4811
$ws.Range("N46").Value = -3659.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H75").Value = 70605
$ws.Range("I75").Value = 30000
$ws.Range("K75").Value = 30000
$ws.Range("M75").Value = -29064

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H78").Value = 70605
$ws.Range("I78").Value = 30000
$ws.Range("K78").Value = 90000
$ws.Range("M78").Value = -85320

$ws.Range("H113").Value = 1415.7333
$ws.Range("I113").Value = 464.30768
$ws.Range("K113").Value = 1392.92304
$ws.Range("M113").Value = 777.0769599999999

$ws.Range("H133").Value = 91000
$ws.Range("J133").Value = 91000
$ws.Range("L133").Value = 91000
$ws.Range("N133").Value = -101120
